$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns remain plain text so Excel does not
# reinterpret values like "1.003" or "314.41" as numbers and round them.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.507.17"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "1.830.98"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "314.41"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").Value = "0.4296"
$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("D8").Value = "0.3659"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "0.07275"
$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("D10").Value = "0.8706"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").Value = "20.66"
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").Value = "1.759.90"
$ws.Range("E12").Value = "  -3.92%  "

$ws.Range("D13").Value = "5.420"
$ws.Range("E13").Value = "  +1.40%  "

$ws.Range("D14").Value = "6.536"
$ws.Range("E14").Value = "  +0.14%  "

$ws.Range("D15").Value = "0.06940"
$ws.Range("E15").Value = "  -0.01%  "

$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").Value = "80.47"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "0.000008934"
$ws.Range("E18").Value = "  -1.15%  "

$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("D20").Value = "15.45"
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").Value = "27.723.73"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").Value = "5.166"
$ws.Range("E22").Value = "  +3.86%  "

$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  +4.72%  "

$ws.Range("D24").Value = "2.100.65"
$ws.Range("E24").Value = "  +3.16%  "

$ws.Range("D25").Value = "1.981"
$ws.Range("E25").Value = "  -0.29%  "

$ws.Range("D26").Value = "154.69"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("D27").Value = "18.86"
$ws.Range("E27").Value = "  +1.12%  "

$ws.Range("D28").Value = "5.168"
$ws.Range("E28").Value = "  -1.68%  "

$ws.Range("D29").Value = "114.40"
$ws.Range("E29").Value = "  -5.39%  "

$ws.Range("E30").Value = "  -1.39%  "

$ws.Range("D31").Value = "0.08905"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").Value = "0.7615"
$ws.Range("E32").Value = "  +0.88%  "

$ws.Range("D33").Value = "4.550"
$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").Value = "2.974"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "1.143"
$ws.Range("E35").Value = "  +1.82%  "

$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.11%  "

$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("D38").Value = "0.05320"
$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("D39").Value = "0.01940"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("D40").Value = "2.803"
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5095"
$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.1667"
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("D43").Value = "6.641"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("D44").Value = "8.426"
$ws.Range("E44").Value = "  +1.24%  "

$ws.Range("D45").Value = "10.52"
$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("E46").Value = "  +1.99%  "

$ws.Range("D47").Value = "0.06509"
$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("D48").Value = "0.4695"
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("D49").Value = "1.002"
$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("D50").Value = "1.620"
$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("D51").Value = "1.757"
$ws.Range("E51").Value = "  +3.58%  "
